# Rename columns "Then_Question" -> "Then_Goto" and "Else_Question" -> "Else_Goto"
# on the active (first) worksheet of the Survey template workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: I1 = Then_Question -> Then_Goto, J1 = Else_Question -> Else_Goto
$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"

# Move/record the active selection to I1, matching the saved workbook state.
$ws.Range("I1").Select()
